$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the "before" data (columns F:V) for every row that is
#     involved in a reshuffle, BEFORE any writes happen. ---
$row61 = $ws.Range("F61:V61").Value()
$row63 = $ws.Range("F63:V63").Value()
$row64 = $ws.Range("F64:V64").Value()

$row67 = $ws.Range("F67:V67").Value()
$row69 = $ws.Range("F69:V69").Value()
$row70 = $ws.Range("F70:V70").Value()

$row72 = $ws.Range("F72:V72").Value()
$row74 = $ws.Range("F74:V74").Value()

$row93 = $ws.Range("F93:V93").Value()
$row94 = $ws.Range("F94:V94").Value()

# --- Re-sort ties: rotate / swap match data among rows that share the
#     same match date, now that two later matches have been inserted
#     into the sheet and shifted the tie-break ordering. ---

# 3-way rotation: 61 <- 63 <- 64 <- 61
$ws.Range("F61:V61").Value = $row63
$ws.Range("F63:V63").Value = $row64
$ws.Range("F64:V64").Value = $row61

# 3-way rotation: 67 <- 69 <- 70 <- 67
$ws.Range("F67:V67").Value = $row69
$ws.Range("F69:V69").Value = $row70
$ws.Range("F70:V70").Value = $row67

# simple swap: 72 <-> 74
$ws.Range("F72:V72").Value = $row74
$ws.Range("F74:V74").Value = $row72

# simple swap: 93 <-> 94
$ws.Range("F93:V93").Value = $row94
$ws.Range("F94:V94").Value = $row93

# --- Append the two newly scraped matches as rows 95 and 96 ---
# Copy row 94's formatting down first (column A / E keep their special
# styles) so the new rows look like every other data row.
$ws.Range("A94:V94").Copy()
$ws.Range("A95:V96").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A95").Value = 94
$ws.Range("B95").Value = "netherlands"
$ws.Range("C95").Value = "tweede-divisie"
$ws.Range("D95").Value = "2023-2024"
$ws.Range("E95").Value = 45241.64583333334
$ws.Range("F95").Value = "Jong Sparta Rotterdam"
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = "AFC"
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 2.38
$ws.Range("K95").Value = "10/11/2023 03:42"
$ws.Range("L95").Value = 2.18
$ws.Range("M95").Value = "11/11/2023 15:16"
$ws.Range("N95").Value = 3.51
$ws.Range("O95").Value = "10/11/2023 03:42"
$ws.Range("P95").Value = 3.57
$ws.Range("Q95").Value = "11/11/2023 15:16"
$ws.Range("R95").Value = 2.37
$ws.Range("S95").Value = "10/11/2023 03:42"
$ws.Range("T95").Value = 2.89
$ws.Range("U95").Value = "11/11/2023 15:16"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/netherlands/tweede-divisie/jong-sparta-rotterdam-afc/OWfoWP0E/"

$ws.Range("A96").Value = 95
$ws.Range("B96").Value = "netherlands"
$ws.Range("C96").Value = "tweede-divisie"
$ws.Range("D96").Value = "2023-2024"
$ws.Range("E96").Value = 45241.75
$ws.Range("F96").Value = "De Treffers"
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = "Hardenberg"
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 1.88
$ws.Range("K96").Value = "10/11/2023 06:12"
$ws.Range("L96").Value = 1.66
$ws.Range("M96").Value = "11/11/2023 17:19"
$ws.Range("N96").Value = 3.6
$ws.Range("O96").Value = "10/11/2023 06:12"
$ws.Range("P96").Value = 3.93
$ws.Range("Q96").Value = "11/11/2023 17:19"
$ws.Range("R96").Value = 3.15
$ws.Range("S96").Value = "10/11/2023 06:12"
$ws.Range("T96").Value = 4.4
$ws.Range("U96").Value = "11/11/2023 17:08"
$ws.Range("V96").Value = "https://www.betexplorer.com/football/netherlands/tweede-divisie/de-treffers-hardenberg/Uy6RBeFa/"
